$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet had a spare row (old row 13, holding only the
# "3268262 - Carlos Renato Menegatti" value in B/C with no A label) that
# was removed, shifting every following row up by one.
$ws.Rows.Item(13).Delete()

# After the shift, several cell contents need to be corrected to match
# the final state of the sheet.
$ws.Range("B10").Value = "3268262 - Carlos Renato Menegatti"
$ws.Range("C10").Value = "3268262 - Carlos Renato Menegatti"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2018" looks like a date to Excel's auto-conversion, so instead
# of assigning it as a literal (which would turn the cell into a date
# serial number and bump its style), copy the already-correctly-typed
# text value from A8/B8 using paste-special values only.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("B18").Value = "3268262 - Carlos Renato Menegatti"
$ws.Range("C18").Value = "3268262 - Carlos Renato Menegatti"

$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
